# "Updated the file location to be static": the productonestatus column
# (C) is re-written cell-by-cell from literal text instead of being left
# as-is, and the column is then re-fit to its (now static/explicit)
# content width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$failedRows = @(11, 17)

for ($row = 2; $row -le 34; $row++) {
    if ($failedRows -contains $row) {
        $ws.Cells.Item($row, 3).Value = "failed"
    } else {
        $ws.Cells.Item($row, 3).Value = "passed"
    }
}

# Column C ("productonestatus") is re-fit to its content width.
$ws.Columns.Item(3).ColumnWidth = 23
